$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (and a Cronos/TrustWalletToken row swap).
# NumberFormat is forced to Text ("@") before each write so numeric-looking
# strings (e.g. "245.99", "0.630", "36.469.84") are preserved verbatim as
# text, matching the original inlineStr cell contents, instead of Excel
# silently coercing them into numbers (which would drop trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.469.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.995.72"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.47"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0775"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.03%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.03"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +13.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.873"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.10"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.301.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.48"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.013.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.393.87"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.65"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.09"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +17.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.06"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.05%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.93"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0615"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.32%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +20.25%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.25"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0992"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.91"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0213"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "92.16"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.56"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.353.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.192.90"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.26%  "
